$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.192674517631531
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.628109812736511
$ws.Range("E1").Value = 1.099163174629211
